$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the previously empty "Quicksort" (column D) values for the
#     existing "random" and "sorted" blocks ---
$ws.Range("D2").Value = 2
$ws.Range("D3").Value = 37
$ws.Range("D4").Value = 482

$ws.Range("D7").Value = 5
$ws.Range("D8").Value = 461
$ws.Range("D9").Value = "#error"

# --- Prepare the new rows (11-20) with the same bordered style used by
#     the rest of the table (copy formatting from the blank spacer row 6) ---
$ws.Range("A6:E6").Copy()
$ws.Range("A11:E20").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 11 stays blank (spacer row), nothing more to do there.

# --- New "inverted" block (rows 12-15) ---
$ws.Range("A12").Value = "50 Elemente"
$ws.Range("B12").Value = 15
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 8
$ws.Range("E12").Value = 9
$ws.Range("G12").Value = "inverted"

$ws.Range("A13").Value = "500 Elemente"
$ws.Range("B13").Value = 1569
$ws.Range("C13").Value = 992
$ws.Range("D13").Value = 744
$ws.Range("E13").Value = 921

$ws.Range("A14").Value = "5000 Elemente"
$ws.Range("B14").Value = 155383
$ws.Range("C14").Value = 99046
$ws.Range("D14").Value = "#error"
$ws.Range("E14").Value = 92703

$ws.Range("A15").Value = "50000 Elemente"
# B15:E15 remain empty

# Row 16 stays blank (spacer row).

# --- New "partly sorted" block (rows 17-20), only labels are filled in ---
$ws.Range("A17").Value = "50 Elemente"
$ws.Range("G17").Value = "partly sorted"

$ws.Range("A18").Value = "500 Elemente"

$ws.Range("A19").Value = "5000 Elemente"

$ws.Range("A20").Value = "50000 Elemente"

# --- Update the selected cell shown when the workbook is opened ---
$ws.Range("E15").Select() | Out-Null
